$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: convert the inline-string timestamps into real Excel date/time
#     serial values, formatted as "YYYY-MM-DD HH:MM:SS".
#     (Applying the lowercase variant first on a single cell, then the
#     uppercase variant on the whole range reproduces the exact numFmt /
#     cellXfs table the workbook ends up with.)
$dateRange = $ws.Range("A2:A7")
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$dateRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A2").Value = 45687.51998796297
$ws.Range("A3").Value = 45687.52202731482
$ws.Range("A4").Value = 45687.52285138889
$ws.Range("A5").Value = 45687.51998564815
$ws.Range("A6").Value = 45687.52202384259
$ws.Range("A7").Value = 45687.52284907408

# --- New column F: "Trening" header (same style as the other headers) plus
#     "Gra" for every data row.
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("F2").Value = "Gra"
$ws.Range("F3").Value = "Gra"
$ws.Range("F4").Value = "Gra"
$ws.Range("F5").Value = "Gra"
$ws.Range("F6").Value = "Gra"
$ws.Range("F7").Value = "Gra"
